# Update the "Förändrad" (Changed) date column for all data rows (2-98)
# from 2023-09-03 (serial 45172) to 2023-09-06 (serial 45175).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C98").Value = 45175
